$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# Title
Replace-Text "Cosmos Unveiled: The Symphony of Existence" "A Journey Through Molecular Medicine: Connecting Health and Chemistry"

# Author name
Replace-Text "Amelia Rose" "Dr. Eleanor Alvarez"

# Email
Replace-Text "cosmicwonder@stardustrealm" "eleanor"
Replace-Text "org" "alvarez@wright.edu"

# Body paragraph 1
Replace-Text "Like a celestial tapestry woven with threads of light, the cosmos presents a mesmerizing spectacle of cosmic harmony" "In the tapestry of human existence, health occupies a central position, entwined with intricate biological processes"
Replace-Text "Each star, planet, and galaxy plays its part in a symphony of existence, an ethereal composition resonating with the echoes of time and space" "Understanding the molecular underpinnings of life offers a key to unravelling the mysteries of disease and paving the path towards novel treatments"
Replace-Text "As we embark on this journey of discovery, we unravel the secrets of the universe, from the depths of black holes to the boundless reaches of nebulae, from the origin of life to the mysteries that lie beyond our mortal grasp" "Molecular medicine, an intersection of chemistry and biology, empowers us to decipher this symphony of life, illuminating the molecular mechanisms that govern health and illness"
Replace-Text "Through the lens of science, philosophy, and art, we seek enlightenment, inspiration, and awe in the grand theater of the cosmos" "As we delve into the realm of molecular medicine, we embark on a journey filled with discoveries, challenges, and boundless opportunities"

# Body paragraph 2
Replace-Text "In this cosmic symphony, celestial bodies dance to the rhythm of gravitational forces, shaping intricate patterns and orchestrating cosmic ballets" "At the heart of molecular medicine lies the chemistry of life"
Replace-Text "The celestial ballet performed by planets as they orbit the Sun, and the graceful pirouette of celestial orbs within galaxies are celestial spectacles of grace and precision" "Molecules, the building blocks of life, engage in intricate interactions, dictating the symphony of biological processes"
Replace-Text "We seek to unravel the mysteries of the unseen forces that sway the celestial symphony, such as dark energy and dark matter, whose enigmatic presence permeates the universe and shapes its ultimate destiny" "From DNA's genetic code to the intricate machinery of proteins, chemistry offers a language to decode the language of life. Molecular medicine empowers us to manipulate these molecules, precisely targeting them to combat disease and promote health. By understanding the molecular mechanisms of disease, we can unravel the enigma of illness, paving the way for effective therapies and interventions"

# Body paragraph 3
Replace-Text "The cosmos hold secrets of our own existence, mirroring the mysteries of life and consciousness within the intricate cosmic script" "Furthermore, molecular medicine offers a unique lens through which we can view and comprehend human health"
Replace-Text "As we explore the nature of space, time, and the fundamental forces of the universe, we search for answers to questions that have haunted humanity for millennia" "It enables us to delve into the molecular interactions of the body, deciphering the intricate dance of cells, tissues, and organs"
Replace-Text "What is the nature of consciousness? Are we alone in the universe? What lies beyond the confines of our visible cosmos? The cosmic journey is not just an exploration of the physical world but a cerebral and spiritual voyage of self-discovery, where we seek to find our place amidst the vastness of the cosmos" "By examining the molecular basis of disease, we gain insights into the complex interplay of genetic, environmental, and lifestyle factors that influence our well-being. This understanding empowers us to intervene at the molecular level, preventing and treating diseases with greater precision and efficacy"

# Summary heading stays the same ("Summary")

# Summary paragraph
Replace-Text "The cosmos, a vast and intricate tapestry of celestial wonders, offers a riveting symphony of existence" "Molecular medicine stands as a beacon of hope in the quest for understanding and treating diseases"
Replace-Text "We embark on a cosmic journey, guided by science, philosophy, and art, to unravel the secrets of the cosmos, from the intricacies of cosmic choreography to the fundamental questions of life and consciousness" "By harnessing the power of chemistry and biology, molecular medicine offers a deeper understanding of the molecular basis of life and illness"
Replace-Text "Through this exploration, we seek enlightenment, awe, and inspiration, and perhaps find our place amidst the grand orchestra of the universe" "It enables us to manipulate molecules, precisely targeting them to combat disease and promote health. Furthermore, molecular medicine provides a unique perspective on human health, allowing us to examine the intricate interactions of the body's molecular machinery. With molecular medicine as our guide, we embark on a transformative journey towards better health and well-being for humanity"

# New trailing empty paragraph
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.Text = "`r"
